$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing runtime values (RuntimesChart, day 1-3 of 2023 part 2) ---
$ws.Range("B3").Value = 0.00848916
$ws.Range("B4").Value = 0.00570892
$ws.Range("B5").Value = 0.024546

# --- Append new days (4-9) with their day number + average runtime ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0.00343778

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.0021247

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 0.01897846

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 0.01155386

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 0.01575866

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 0.01636104

# --- Match the author's final selection over the whole data block ---
$ws.Range("A3:B11").Select() | Out-Null

# --- Extend the bar chart's source ranges to cover the new rows ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(RuntimesChart!`$B`$2,RuntimesChart!`$A`$3:`$A`$13,RuntimesChart!`$B`$3:`$B`$13,1)"
